$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 109 - fill in existing blank template row
$ws.Range("A109").Value = 42086
$ws.Range("B109").Value = 4
$ws.Range("C109").Value = "Vergadering"

# Row 110 - date/hours left blank, only description added
$ws.Range("C110").Value = "zoeken voor oplossingen bugs"

# Row 111 - fill in existing blank template row; A111 also needs its
# style fixed from s=3 (stray) to s=2 (the standard date-cell style)
$ws.Range("A109").Copy() | Out-Null
$ws.Range("A111").PasteSpecial(-4122) | Out-Null
$ws.Range("A111").Value = 42088
$ws.Range("B111").Value = 4
$ws.Range("C111").Value = "subscriber toevoegen"

# Row 112 - new row, copy formatting (s=2,3,4) from row 109, then set values
$ws.Range("A109:C109").Copy() | Out-Null
$ws.Range("A112:C112").PasteSpecial(-4122) | Out-Null
$ws.Range("A112").Value = 42089
$ws.Range("C112").Value = "UX breed in organism"

# Rows 113-116 - new blank rows with same formatting pattern (s=2,3,4) as row 109
$ws.Range("A109:C109").Copy() | Out-Null
$ws.Range("A113:C116").PasteSpecial(-4122) | Out-Null

# Row 117 - new blank row with formatting pattern (s=3,3,4), matching row 4's pattern
$ws.Range("A4:C4").Copy() | Out-Null
$ws.Range("A117:C117").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Update sheet view to reflect scrolled position / active selection
$ws.Application.ActiveWindow.ScrollRow = 79
$ws.Range("C112").Select() | Out-Null
